$p = $ppt.ActivePresentation

# Reorder slides.
# Current order (1-indexed): 1 Kiddo, 2 What is Kiddo?, 3 Why Kiddo?,
#   4 Who?(Wanda), 5 Who?(Ted), 6 Who?(Pat), 7 What we have completed so far,
#   8 What we plan to do next
# Target order: 1 Kiddo, 2 What is Kiddo?, 3 Who?(Wanda), 4 Who?(Ted),
#   5 Who?(Pat), 6 Why Kiddo?
# Move "Why Kiddo?" (slide 3) to the end of the deck.
$p.Slides.Item(3).MoveTo($p.Slides.Count)

# Delete the trailing "next steps" slides (now at positions 6 and 7,
# since "Why Kiddo?" was appended as the new last slide).
$p.Slides.Item(7).Delete()
$p.Slides.Item(6).Delete()
